$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - columns C and D swap meaning:
#   C1 was "appt time", now becomes "doc location"
#   D1 was "doc location", now becomes "apt time" (re-spelled)
$ws.Range("C1").Value = "doc location"
$ws.Range("D1").Value = "apt time"

# Row 2
$ws.Range("A2").Value = "812-121-0912"
$ws.Range("B2").Value = "Dr. Bailey"
$ws.Range("C2").Value = "Watson Clinic South"
$ws.Range("D2").Value = 0.38541666666666669
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("E2").Value = $false
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = $true
$ws.Range("H2").Value = 0

# Row 3
$ws.Range("A3").Value = "503-388-1908"
$ws.Range("B3").Value = "Dr. Miley"
$ws.Range("C3").Value = "East Bay Clinic"
$ws.Range("D3").Value = 0.48958333333333331
$ws.Range("D3").NumberFormat = "h:mm"
$ws.Range("E3").Value = $true
$ws.Range("F3").Value = 125.75
$ws.Range("G3").Value = $true
$ws.Range("H3").Value = 49.54

# Row 4
$ws.Range("A4").Value = "810-225-7205"
$ws.Range("B4").Value = "Dr. Hurtak"
$ws.Range("C4").Value = "Tampa General"
$ws.Range("D4").Value = 0.33333333333333331
$ws.Range("D4").NumberFormat = "h:mm"
$ws.Range("E4").Value = $true
$ws.Range("F4").Value = 100.34
$ws.Range("G4").Value = $true
$ws.Range("H4").Value = 75.12

# Row 5
$ws.Range("A5").Value = "011-433-3770"
$ws.Range("B5").Value = "Dr.Bailey"
$ws.Range("C5").Value = "Watson Clinic South"
$ws.Range("D5").Value = 0.46875
$ws.Range("D5").NumberFormat = "h:mm"
$ws.Range("E5").Value = $false
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = $false
$ws.Range("H5").Value = 0

# Row 6
$ws.Range("A6").Value = "800-992-2131"
$ws.Range("B6").Value = "Dr. Mallove"
$ws.Range("C6").Value = "Southshore Reginal"
$ws.Range("D6").Value = 0.10416666666666667
$ws.Range("D6").NumberFormat = "h:mm"
$ws.Range("E6").Value = $true
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = $true
$ws.Range("H6").Value = 50

$ws.Range("F6").Select()
